# Applies the crypto price/volume/name refresh described by the commit:
# "Updated cryptos list on Sun Feb 25 17:57:55 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.544.03"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.062.16"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "385.88"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.35"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  -0.57%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.92"
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0862"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.541.32"
$ws.Range("E13").Value = "  +2.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.58"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.78"
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.060.03"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.976"
$ws.Range("E17").Value = "  -2.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.73"
$ws.Range("E18").Value = "  -4.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.639.23"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("E21").Value = "  -1.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0965"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.24"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.17"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.25"
$ws.Range("E26").Value = "  +5.26%  "
$ws.Range("E27").Value = "  +2.73%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.28"
$ws.Range("E28").Value = "  -1.85%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.170"
$ws.Range("E29").Value = "  +2.33%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.28"
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.72"
$ws.Range("E33").Value = "  -0.63%  "
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.44"
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0447"
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("E38").Value = "  +2.37%  "
$ws.Range("E39").Value = "  +7.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.00"
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.56"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.31"
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("E45").Value = "  +2.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.94"
$ws.Range("E46").Value = "  +1.53%  "
$ws.Range("E47").Value = "  +3.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.43"
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.033.76"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.355.83"
$ws.Range("E50").Value = "  +2.08%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.206"
$ws.Range("E51").Value = "  +5.93%  "
